$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 474 (a new weekly reading, dated 2021-11-05 / serial 44505),
# pre-populating them as a copy of the row block that will now follow them (old row 474-476,
# now shifted to 477-479) so that formatting/structure is consistent, then overwrite the
# values that differ for the new week.
$ws.Rows("474:476").Insert()
$ws.Range("A474:R476").Value2 = $ws.Range("A477:R479").Value2

# New values for the inserted rows (date 44505 = 2021-11-05)
$ws.Range("D474").Value2 = 44505
$ws.Range("J474").Value2 = 48000
$ws.Range("K474").Value2 = 90
$ws.Range("L474").Value2 = 100
$ws.Range("M474").Value2 = 95
$ws.Range("P474").Value2 = 95

$ws.Range("D475").Value2 = 44505
$ws.Range("J475").Value2 = 40000
$ws.Range("K475").Value2 = 75
$ws.Range("L475").Value2 = 80
$ws.Range("M475").Value2 = 77
$ws.Range("P475").Value2 = 77

$ws.Range("D476").Value2 = 44505
$ws.Range("J476").Value2 = 15000
$ws.Range("K476").Value2 = 60
$ws.Range("L476").Value2 = 60
$ws.Range("M476").Value2 = 60
$ws.Range("P476").Value2 = 60

Write-Output "done"
